$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: restyle row 27 (A27 gains a bottom-border cell; B-E27 switch to the
#     bottom-border variant of their current style) by copying formats from row 5,
#     which already carries that exact bottom-border style. ---
$ws.Range("A5:E5").Copy() | Out-Null
$ws.Range("A27:E27").PasteSpecial(-4122) | Out-Null

# --- Step 2: pre-format the 6 new rows (28-33) by copying the plain (no-border)
#     style from row 3 - a ht=43.2, 3-line template matching the new rows. ---
$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A28:E33").PasteSpecial(-4122) | Out-Null

# --- Step 3: row 31 is the last row of its filename-group, so it additionally gets
#     the bottom-border style (same donor as row 27). ---
$ws.Range("A5:E5").Copy() | Out-Null
$ws.Range("A31:E31").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Step 4: row heights for the new rows (3 wrapped lines each) ---
foreach ($r in 28..33) {
  $ws.Rows.Item($r).RowHeight = 43.2
}

# --- Step 5: cell values for the new rows ---
$ws.Cells.Item(28,1).Value = 'SCRIPT/D73P31A/us3103.ssb '
$ws.Cells.Item(28,2).Value = 54
$ws.Cells.Item(28,3).Value = ' Yay! Hooray!\nWe''ve reached the summit!'
$ws.Cells.Item(28,4).Value = ' Да! Урааа! Мы достигли вершины!'
$ws.Cells.Item(28,5).Value = ' Äà! Ôñààà! Íú äïòóéãìé âåñšéîú!'

$ws.Cells.Item(29,1).Value = 'SCRIPT/D73P31A/us0202.ssb'
$ws.Cells.Item(29,2).Value = 57
$ws.Cells.Item(29,3).Value = ' I joined up with [CS:N]Heracross[CR] at the\n7th Station Clearing, and we were climbing\ntogether. It was such fun!'
$ws.Cells.Item(29,4).Value = ' Я встретил [CS:N]Геракросса[CR] на Поляне\n7-го Перехода и мы пошли вместе. Это было\nтак весело!'
$ws.Cells.Item(29,5).Value = ' Ÿ âòóñåóéì [CS:N]Ãåñàëñïòòà[CR] îà Ðïìÿîå\n7-ãï Ðåñåöïäà é íú ðïšìé âíåòóå. Üóï áúìï\nóàë âåòåìï!'

$ws.Cells.Item(30,1).Value = 'SCRIPT/D73P31A/us0302.ssb'
$ws.Cells.Item(30,2).Value = 60
$ws.Cells.Item(30,3).Value = ' He kept yelling, \"Whooooa!\"\nJust watching it was so much fun.'
$ws.Cells.Item(30,4).Value = ' Он постоянно орал: \"Ваааааау!\"\nМне было так весело это наблюдать.'
$ws.Cells.Item(30,5).Value = ' Ïî ðïòóïÿîîï ïñàì: \"Âààààààô!\"\nÍîå áúìï óàë âåòåìï üóï îàáìýäàóû.'

$ws.Cells.Item(31,1).Value = 'SCRIPT/D73P31A/us0402.ssb'
$ws.Cells.Item(31,2).Value = 63
$ws.Cells.Item(31,3).Value = ' I''ve heard it so many times,\nI can imitate it. Ha ha ha!'
$ws.Cells.Item(31,4).Value = ' Я слышал этот крик столько раз,\nчто смогу в точности повторить его.\nХа-ха-ха!'
$ws.Cells.Item(31,5).Value = ' Ÿ òìúšàì üóïó ëñéë òóïìûëï ñàè,\nœóï òíïãô â óïœîïòóé ðïâóïñéóû åãï.\nÖà-öà-öà!'

$ws.Cells.Item(32,1).Value = 'SCRIPT/P01P04A/us2004.ssb '
$ws.Cells.Item(32,2).Value = 38
$ws.Cells.Item(32,3).Value = ' [CS:P]Sky Peak[CR] was so much fun…'
$ws.Cells.Item(32,4).Value = ' Мне так понравился [CS:P]Небесный\nПик[CR]...'
$ws.Cells.Item(32,5).Value = ' Íîå óàë ðïîñàâéìòÿ [CS:P]Îåáåòîúê\nÐéë[CR]...'

$ws.Cells.Item(33,1).Value = 'SCRIPT/P01P04A/us2104.ssb '
$ws.Cells.Item(33,2).Value = 41
$ws.Cells.Item(33,3).Value = ' I wonder where I''ll explore\nnext time.'
$ws.Cells.Item(33,4).Value = ' Интересно, какое место я буду\nисследовать в следующий раз.'
$ws.Cells.Item(33,5).Value = ' Éîóåñåòîï, ëàëïå íåòóï ÿ áôäô\néòòìåäïâàóû â òìåäôýþéê ñàè.'

# --- Step 6: viewport / selection, matching the author's final cursor position ---
$ws.Range("E33").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 31
